$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.722.27"
$ws.Range("E2").Value = "  -2.86%  "
$ws.Range("D3").Value = "'3.440.87"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'582.88"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").Value = "'172.61"
$ws.Range("E6").Value = "  -3.84%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("D9").Value = "'3.437.96"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").Value = "'0.130"
$ws.Range("E10").Value = "  -5.92%  "
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").Value = "'0.408"
$ws.Range("E12").Value = "  -4.51%  "
$ws.Range("D13").Value = "'4.035.97"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "'28.74"
$ws.Range("E15").Value = "  -11.61%  "
$ws.Range("D16").Value = "'65.808.86"
$ws.Range("E16").Value = "  -2.68%  "
$ws.Range("D17").Value = "'0.0000170"
$ws.Range("E17").Value = "  -3.91%  "
$ws.Range("D18").Value = "'3.436.26"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").Value = "'5.92"
$ws.Range("E19").Value = "  -3.78%  "
$ws.Range("D20").Value = "'13.83"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Value = "'365.64"
$ws.Range("E21").Value = "  -6.42%  "
$ws.Range("D22").Value = "'7.66"
$ws.Range("D23").Value = "'72.73"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'0.533"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").Value = "'0.0000122"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").Value = "'9.76"
$ws.Range("E27").Value = "  -4.15%  "
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'23.72"
$ws.Range("E30").Value = "  -4.29%  "
$ws.Range("D31").Value = "'1.98"
$ws.Range("E31").Value = "  -3.34%  "
$ws.Range("D32").Value = "'5.69"
$ws.Range("E32").Value = "  -6.46%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -7.05%  "
$ws.Range("D35").Value = "'7.02"
$ws.Range("E35").Value = "  -3.78%  "
$ws.Range("D36").Value = "'1.52"
$ws.Range("E36").Value = "  -2.83%  "
$ws.Range("D37").Value = "'161.37"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "'29.02"
$ws.Range("E38").Value = "  +3.20%  "
$ws.Range("D39").Value = "'0.881"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").Value = "'2.60"
$ws.Range("E40").Value = "  -3.73%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.75"
$ws.Range("E41").Value = "  -5.71%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'2.760.18"
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("D43").Value = "'4.43"
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("D44").Value = "'6.42"
$ws.Range("D45").Value = "'0.0680"
$ws.Range("E45").Value = "  -4.89%  "
$ws.Range("D46").Value = "'40.00"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("D47").Value = "'24.11"
$ws.Range("D48").Value = "'0.0289"
$ws.Range("E48").Value = "  -3.19%  "
$ws.Range("D49").Value = "'321.72"
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("D51").Value = "'6.20"
$ws.Range("E51").Value = "  -1.36%  "
